$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 9.6

# Row 4
$ws.Range("G4").Value = 5.5
$ws.Range("H4").Value = 3.55
$ws.Range("I4").Value = 1.6
$ws.Range("J4").Value = 5.6
$ws.Range("K4").Value = 2.12
$ws.Range("L4").Value = 2.12
$ws.Range("O4").Value = 1.37
$ws.Range("Q4").Value = 2.07
$ws.Range("T4").Value = 2.52
$ws.Range("U4").Value = 2.07
$ws.Range("V4").Value = 1.6
$ws.Range("W4").Value = 12
$ws.Range("X4").Value = 32
$ws.Range("Y4").Value = 18.5
$ws.Range("Z4").Value = 110
$ws.Range("AA4").Value = 70
$ws.Range("AB4").Value = 75
$ws.Range("AD4").Value = 7.1
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 120
$ws.Range("AH4").Value = 5.4
$ws.Range("AI4").Value = 6.5
$ws.Range("AJ4").Value = 8.5
$ws.Range("AK4").Value = 11.25
$ws.Range("AL4").Value = 14.5
$ws.Range("AM4").Value = 35
$ws.Range("AN4").Value = 6.9
$ws.Range("AO4").Value = 32
$ws.Range("AP4").Value = 40
$ws.Range("AU4").Value = 8
$ws.Range("AV4").Value = 90
$ws.Range("AW4").Value = 3.25
$ws.Range("AX4").Value = 7.5
$ws.Range("AY4").Value = 19
$ws.Range("AZ4").Value = 25
$ws.Range("BB4").Value = 300

# Row 5
$ws.Range("G5").Value = 3.1
$ws.Range("I5").Value = 2.2
$ws.Range("J5").Value = 3.6
$ws.Range("U5").Value = 1.67
$ws.Range("V5").Value = 2.1
$ws.Range("Z5").Value = 34
$ws.Range("AE5").Value = 13
$ws.Range("AH5").Value = 9
$ws.Range("AO5").Value = 17
$ws.Range("AW5").Value = 4.33
